$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62: the B/F/G cells were blank placeholders; drop them entirely ---
$ws.Cells.Item(62, 2).ClearContents()   # B62
$ws.Cells.Item(62, 6).ClearContents()   # F62
$ws.Cells.Item(62, 7).ClearContents()   # G62

# --- Row 63 ---
$ws.Cells.Item(63, 1).Value = '05/01/2026 10:30:14'
$ws.Cells.Item(63, 2).Value = '05/01 10:11'
$ws.Cells.Item(63, 3).Value = 'Metrópoles'
$ws.Cells.Item(63, 4).Value = 'Empresa de filho de Lula deve mais de R$ 370 mil à União'
$ws.Cells.Item(63, 5).Value = 'https://www.metropoles.com/colunas/andreza-matais/empresa-de-filho-de-lula-deve-mais-de-r-370-mil-a-uniao'
$ws.Cells.Item(63, 6).Value = 'tribut'
$ws.Cells.Item(63, 7).Value = 'Revendedora de gás de Sandro Luís Lula da Silva é alvo de execuções judiciais por débitos &lt;b&gt;tribut&lt;/b&gt;ários e previdenciários'

# --- Row 64 ---
$ws.Cells.Item(64, 1).Value = '05/01/2026 10:30:15'
$ws.Cells.Item(64, 2).Value = '05/01 10:08'
$ws.Cells.Item(64, 3).Value = 'Metrópoles'
$ws.Cells.Item(64, 4).Value = 'CGU identifica desvios e falta de controle em viagens na Codeba'
$ws.Cells.Item(64, 5).Value = 'https://www.metropoles.com/colunas/tacio-lorran/cgu-identifica-desvios-e-falta-de-controle-em-viagens-na-codeba'
$ws.Cells.Item(64, 6).Value = 'cgu'
$ws.Cells.Item(64, 7).Value = 'Denúncia sobre gastos excessivos de diretores da Codeba em viagens motivou escrutínio da CGU'

# --- Row 65 ---
$ws.Cells.Item(65, 1).Value = '05/01/2026 10:30:15'
$ws.Cells.Item(65, 2).Value = '05/01 10:00'
$ws.Cells.Item(65, 3).Value = 'Metrópoles'
$ws.Cells.Item(65, 4).Value = '"Inferno Catarina”: pré-candidato promete “pior prisão do país" em SC'
$ws.Cells.Item(65, 5).Value = 'https://www.metropoles.com/colunas/paulo-cappelli/pre-candidato-mestre-de-luta-livre-propoe-presidio-inferno-catarina'
$ws.Cells.Item(65, 6).Value = 'santa catarina'
$ws.Cells.Item(65, 7).Value = 'Pré-candidato em Santa Catarina, o mestre de Luta Livre Marcelo Brigadeiro quer "tornar a vida dos detentos a mais sofrid'

# --- Row 66 ---
$ws.Cells.Item(66, 1).Value = '05/01/2026 10:30:16'
$ws.Cells.Item(66, 2).Value = '05/01 10:00'
$ws.Cells.Item(66, 3).Value = 'Folha de S.Paulo - Poder - Principal'
$ws.Cells.Item(66, 4).Value = 'Governo Castro prevê novas fases de operação que deixou 122 mortos'
$ws.Cells.Item(66, 5).Value = 'https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.com.br/colunas/painel/2026/01/governo-castro-preve-novas-fases-de-operacao-que-deixou-122-mortos.shtml'
$ws.Cells.Item(66, 6).Value = 'senado'
# G66 was captured upstream as a broken/truncated formula (source feed text that
# happened to start with a quote) rather than a plain string - reproduce it as a
# formula so the cell keeps the same <f> shape as the source workbook.
$ws.Cells.Item(66, 7).Formula = "=`"https://www1.folha.uol.com.br/folha-topicos/pl/`"&gt;PL&lt;/a&gt;) deve renunciar para disputar o Senado.
&lt;a href=`"https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.co"

# --- Row 67 ---
$ws.Cells.Item(67, 1).Value = '05/01/2026 10:30:17'
# B67 is an explicit-but-blank field in the source feed (column present, no value)
$ws.Cells.Item(67, 2).Style = "Normal"
$ws.Cells.Item(67, 3).Value = 'VEJA'
$ws.Cells.Item(67, 4).Value = 'Quem é o secretário de Haddad que deixou o governo'
$ws.Cells.Item(67, 5).Value = 'https://veja.abril.com.br/coluna/radar-economico/quem-e-o-secretario-do-haddad-que-deixou-o-governo/'
$ws.Cells.Item(67, 6).Value = 'ministério da fazenda'
$ws.Cells.Item(67, 7).Value = 'bosa Pinto deixou nesta segunda-feira, 5, o cargo de secretário de Reformas Econômicas do Ministério da Fazenda'

Write-Host "applied historico_news update"
